$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3844.5454
$ws.Range("I76").Value = 3249.5
$ws.Range("K76").Value = 3249.5
$ws.Range("M76").Value = -2934.5
$ws.Range("H79").Value = 3844.5454
$ws.Range("I79").Value = 3249.5
$ws.Range("K79").Value = 3249.5
$ws.Range("M79").Value = -2157.5
$ws.Range("H80").Value = 748.6
$ws.Range("I80").Value = 781.5
$ws.Range("J80").Value = 711
$ws.Range("K80").Value = 2344.5
$ws.Range("L80").Value = 2133
$ws.Range("M80").Value = -1346.5
$ws.Range("N80").Value = -4129
$ws.Range("H83").Value = 748.6
$ws.Range("I83").Value = 781.5
$ws.Range("J83").Value = 711
$ws.Range("K83").Value = 7033.5
$ws.Range("L83").Value = 6399
$ws.Range("M83").Value = -2041.5
$ws.Range("N83").Value = -16383
$ws.Range("H88").Value = 1225
$ws.Range("I88").Value = 1997.5
$ws.Range("J88").Value = 916
$ws.Range("K88").Value = 1997.5
$ws.Range("L88").Value = 916
$ws.Range("M88").Value = -1591.5
$ws.Range("N88").Value = -1728
$ws.Range("H91").Value = 1225
$ws.Range("I91").Value = 1997.5
$ws.Range("J91").Value = 916
$ws.Range("K91").Value = 1997.5
$ws.Range("L91").Value = 916
$ws.Range("M91").Value = -593.5
$ws.Range("N91").Value = -3724
$ws.Range("H111").Value = 50000200
$ws.Range("I111").Value = 50000200
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 150000600
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -149997533
$ws.Range("N111").ClearContents()
$ws.Range("H129").Value = 927.7436
$ws.Range("J129").Value = 915.0571
$ws.Range("L129").Value = 2745.1713
$ws.Range("N129").Value = -12745.1713
$ws.Range("H137").Value = 53712.684
$ws.Range("I137").Value = 692.75
$ws.Range("J137").Value = 144604
$ws.Range("K137").Value = 2078.25
$ws.Range("L137").Value = 433812
$ws.Range("M137").Value = 471.75
$ws.Range("N137").Value = -438912

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5000025
$ws.Range("I6").Value = 5000025
$ws.Range("K6").Value = 5000025
$ws.Range("M6").Value = -4999852
$ws.Range("H45").Value = 999
$ws.Range("I45").Value = 999
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 999
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -622
$ws.Range("N45").ClearContents()
$ws.Range("H122").Value = 1736.0588
$ws.Range("I122").Value = 1769.5625
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 5308.6875
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -2858.6875
$ws.Range("N122").Value = -8500
$ws.Range("H132").Value = 2919.3333
$ws.Range("I132").Value = 2965.818
$ws.Range("K132").Value = 8897.454000000002
$ws.Range("M132").Value = -6367.454000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 20274.727
$ws.Range("J11").Value = 36602.5
$ws.Range("L11").Value = 36602.5
$ws.Range("N11").Value = -36882.5
$ws.Range("H16").Value = 9500.5
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 9500.5
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 9500.5
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -9840.5
$ws.Range("H19").Value = 7777.6665
$ws.Range("J19").Value = 9999.799999999999
$ws.Range("L19").Value = 9999.799999999999
$ws.Range("N19").Value = -10345.8
$ws.Range("H87").Value = 25000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 25000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -887
$ws.Range("H31").Value = 2097.535
$ws.Range("J31").Value = 2760.52
$ws.Range("L31").Value = 2760.52
$ws.Range("N31").Value = -3350.52
$ws.Range("H34").Value = 2097.535
$ws.Range("J34").Value = 2760.52
$ws.Range("L34").Value = 2760.52
$ws.Range("N34").Value = -3164.52

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 54.2
$ws.Range("I23").Value = 80
$ws.Range("J23").Value = 47.75
$ws.Range("K23").Value = 240
$ws.Range("L23").Value = 143.25
$ws.Range("M23").Value = -5
$ws.Range("N23").Value = -613.25
$ws.Range("H39").Value = 2344.6365
$ws.Range("J39").Value = 2604.4736
$ws.Range("L39").Value = 7813.4208
$ws.Range("N39").Value = -8401.4208
$ws.Range("H61").Value = 228.57143
$ws.Range("J61").Value = 240
$ws.Range("L61").Value = 720
$ws.Range("N61").Value = -1150
$ws.Range("H100").Value = 3402.1667
$ws.Range("J100").Value = 4222
$ws.Range("L100").Value = 12666
$ws.Range("N100").Value = -14288
$ws.Range("H131").Value = 12518697
$ws.Range("J131").Value = 19654.342
$ws.Range("L131").Value = 58963.026
$ws.Range("N131").Value = -69043.026

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 24999.666
$ws.Range("J49").Value = 24999.666
$ws.Range("L49").Value = 24999.666
$ws.Range("N49").Value = -25367.666
$ws.Range("H54").Value = 19500
$ws.Range("J54").Value = 19500
$ws.Range("L54").Value = 19500
$ws.Range("N54").Value = -20280
$ws.Range("H70").Value = 20400
$ws.Range("H73").Value = 20400
$ws.Range("H132").Value = 1101706.2
$ws.Range("I132").Value = 1925053.2
$ws.Range("J132").Value = 3910.0667
$ws.Range("K132").Value = 5775159.6
$ws.Range("L132").Value = 11730.2001
$ws.Range("M132").Value = -5772629.6
$ws.Range("N132").Value = -16790.2001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 24155.166
$ws.Range("I42").Value = 22000
$ws.Range("J42").Value = 24586.2
$ws.Range("K42").Value = 22000
$ws.Range("L42").Value = 24586.2
$ws.Range("M42").Value = -21437
$ws.Range("N42").Value = -25712.2
$ws.Range("H46").Value = 970.5238000000001
$ws.Range("I46").Value = 486.875
$ws.Range("J46").Value = 1268.1538
$ws.Range("K46").Value = 486.875
$ws.Range("L46").Value = 1268.1538
$ws.Range("M46").Value = -298.875
$ws.Range("N46").Value = -1644.1538
$ws.Range("H49").Value = 24155.166
$ws.Range("I49").Value = 22000
$ws.Range("J49").Value = 24586.2
$ws.Range("K49").Value = 22000
$ws.Range("L49").Value = 24586.2
$ws.Range("M49").Value = -21853
$ws.Range("N49").Value = -24880.2
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H132").Value = 3375.68
$ws.Range("I132").Value = 1300
$ws.Range("J132").Value = 5624.3335
$ws.Range("K132").Value = 3900
$ws.Range("L132").Value = 16873.0005
$ws.Range("M132").Value = -1370
$ws.Range("N132").Value = -21933.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27496
$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -87480
